# Update imputed KNN result values for the terrestrial_mammals / ACD / seed4 dataset.
# Only the numeric values listed below changed between the two versions of the
# workbook; everything else (headers, formatting, other cells) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.045
$ws.Range("A9").Value = -21.831
$ws.Range("C12").Value = -11.277
$ws.Range("D15").Value = -8.348000000000001
$ws.Range("A18").Value = -22.247
$ws.Range("A20").Value = -20.295
$ws.Range("C26").Value = -13.131
$ws.Range("A27").Value = -21.761
$ws.Range("C27").Value = -13.363
$ws.Range("C29").Value = -12.124
$ws.Range("C37").Value = -13.351
$ws.Range("C38").Value = -13.818
$ws.Range("D38").Value = -7.935
$ws.Range("D44").Value = -7.417999999999999
$ws.Range("C51").Value = -12.405
$ws.Range("D51").Value = -7.556
$ws.Range("C55").Value = -13.752
$ws.Range("D57").Value = -8.032
$ws.Range("D63").Value = -7.337000000000001
$ws.Range("A69").Value = -21.831
$ws.Range("C69").Value = -11.627
$ws.Range("C70").Value = -12.716
$ws.Range("D70").Value = -7.858
$ws.Range("A76").Value = -20.306
$ws.Range("A82").Value = -22.205
$ws.Range("C83").Value = -13.551
$ws.Range("D99").Value = -8.103999999999999
$ws.Range("C102").Value = -13.419
